# "mise à jour de la date" -- bump the cached date-field text on the
# master/layouts from 25/02/2019 to 26/02/2019, and fix up the title
# slide: merge the title runs into one, and split the "xx/03/2019"
# placeholder into "06" + "/03/2019".

$p = $ppt.ActivePresentation

function Update-DateField($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq "25/02/2019") {
                $tr.Text = "26/02/2019"
            }
        }
    }
}

# 1) The slide master's "Date Placeholder" auto-date field.
$master = $p.SlideMaster
Update-DateField $master.Shapes

# 2) Same field on every slide layout inheriting from the master.
for ($j = 1; $j -le $master.CustomLayouts.Count; $j++) {
    $layout = $master.CustomLayouts.Item($j)
    Update-DateField $layout.Shapes
}

# 3) Title slide touch-ups.
$slide1 = $p.Slides.Item(1)

# 3a) "Soutenance Projet " + "6" -> single run "Soutenance Projet 6".
$titleRange = $slide1.Shapes.Item(1).TextFrame.TextRange
if ($titleRange.Text -eq "Soutenance Projet 6") {
    $titleRange.Text = "~"
    $titleRange.Text = "Soutenance Projet 6"
}

# 3b) "xx/03/2019" -> "06" + "/03/2019" (two runs).
$subtitleRange = $slide1.Shapes.Item(2).TextFrame.TextRange
$datePara = $subtitleRange.Paragraphs(2, 1)
if ($datePara.Text -eq "xx/03/2019") {
    $datePrefix = $datePara.Characters(1, 2)
    $datePrefix.Text = "06"
}

Write-Output "done"
